$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14
$ws.Range("E2").Value = 14
$ws.Range("F2").Value = 13
$ws.Range("AD2").Value = 2065
$ws.Range("AE2").Value = 1.31
$ws.Range("AG2").Value = 0.47
$ws.Range("AH2").Value = 1646
$ws.Range("AL2").Value = 47619
$ws.Range("AM2").Value = 0.27
$ws.Range("AO2").Value = 0.98
$ws.Range("AP2").Value = 3472
$ws.Range("AX2").Value = 0.84
$ws.Range("AZ2").Value = 41059
$ws.Range("BA2").Value = 52
$ws.Range("BB2").Value = 437
$ws.Range("BC2").Value = 817
$ws.Range("BD2").Value = 966
$ws.Range("BE2").Value = 1546
$ws.Range("Q3").Value = 1457
$ws.Range("Y3").Value = "84"
$ws.Range("Z3").Value = "1541"
$ws.Range("AA3").Value = "254 x 102 x 22"
$ws.Range("AC3").Value = 1988
$ws.Range("AD3").Value = 1877
$ws.Range("AE3").Value = 1.03
$ws.Range("AG3").Value = 0.65
$ws.Range("AH3").Value = 1283
$ws.Range("AJ3").Value = 1988
$ws.Range("AK3").Value = 140000
$ws.Range("AL3").Value = 92528
$ws.Range("AP3").Value = 2011
$ws.Range("AR3").Value = 994
$ws.Range("AS3").Value = 3854
$ws.Range("AT3").Value = 0.51
$ws.Range("AV3").Value = 0.92
$ws.Range("AW3").Value = 916
$ws.Range("AX3").Value = 0.67
$ws.Range("AZ3").Value = 92528
$ws.Range("BD3").Value = 227
$ws.Range("BE3").Value = 793
$ws.Range("B4").Value = 14
$ws.Range("E4").Value = 14
$ws.Range("F4").Value = 13
$ws.Range("AA4").Value = "457 x 152 x 60"
$ws.Range("AC4").Value = 5410
$ws.Range("AD4").Value = 4213
$ws.Range("AE4").Value = 1.13
$ws.Range("AG4").Value = 0.57
$ws.Range("AH4").Value = 3103
$ws.Range("AJ4").Value = 5410
$ws.Range("AK4").Value = 381000
$ws.Range("AL4").Value = 73009
$ws.Range("AM4").Value = 0.27
$ws.Range("AO4").Value = 0.98
$ws.Range("AP4").Value = 5323
$ws.Range("AR4").Value = 2705
$ws.Range("AS4").Value = 25746
$ws.Range("AT4").Value = 0.32
$ws.Range("AV4").Value = 0.97
$ws.Range("AW4").Value = 2629
$ws.Range("AX4").Value = 1.04
$ws.Range("AZ4").Value = 62952
$ws.Range("BA4").Value = 52
$ws.Range("BB4").Value = 437
$ws.Range("BC4").Value = 817
$ws.Range("BD4").Value = 1023
$ws.Range("BE4").Value = 2431
$ws.Range("Q5").Value = 2470
$ws.Range("Z5").Value = "2470"
$ws.Range("AA5").Value = "203 x 133 x 30"
$ws.Range("AC5").Value = 2712
$ws.Range("AD5").Value = 1917
$ws.Range("AE5").Value = 1.19
$ws.Range("AG5").Value = 0.54
$ws.Range("AH5").Value = 1455
$ws.Range("AJ5").Value = 2712
$ws.Range("AK5").Value = 191000
$ws.Range("AL5").Value = 126234
$ws.Range("AP5").Value = 2744
$ws.Range("AR5").Value = 1356
$ws.Range("AS5").Value = 12468
$ws.Range("AT5").Value = 0.33
$ws.Range("AV5").Value = 0.97
$ws.Range("AW5").Value = 1316
$ws.Range("AX5").Value = 0.75
$ws.Range("AZ5").Value = 126234
$ws.Range("BD5").Value = 241
$ws.Range("BE5").Value = 1259
$ws.Range("B6").Value = 14
$ws.Range("E6").Value = 14
$ws.Range("F6").Value = 13
$ws.Range("K6").Value = 1510
$ws.Range("P6").Value = 9.630000000000001
$ws.Range("AD6").Value = 10674
$ws.Range("AE6").Value = 0.88
$ws.Range("AG6").Value = 0.75
$ws.Range("AH6").Value = 6193
$ws.Range("AL6").Value = 112101
$ws.Range("AM6").Value = 0.27
$ws.Range("AO6").Value = 0.98
$ws.Range("AP6").Value = 8173
$ws.Range("AX6").Value = 1.36
$ws.Range("AY6").Value = 224223560
$ws.Range("AZ6").Value = 96658
$ws.Range("BA6").Value = 52
$ws.Range("BB6").Value = 437
$ws.Range("BC6").Value = 817
$ws.Range("BD6").Value = 1128
$ws.Range("BE6").Value = 4037
$ws.Range("K7").Value = 1510
$ws.Range("P7").Value = 9.630000000000001
$ws.Range("Q7").Value = 4162
$ws.Range("Z7").Value = "4162"
$ws.Range("AA7").Value = "356 x 171 x 51"
$ws.Range("AC7").Value = 4608
$ws.Range("AD7").Value = 9319
$ws.Range("AE7").Value = 0.7
$ws.Range("AG7").Value = 0.85
$ws.Range("AH7").Value = 3899
$ws.Range("AJ7").Value = 4608
$ws.Range("AK7").Value = 324500
$ws.Range("AL7").Value = 214466
$ws.Range("AP7").Value = 4661
$ws.Range("AR7").Value = 2304
$ws.Range("AS7").Value = 31348
$ws.Range("AT7").Value = 0.27
$ws.Range("AW7").Value = 2267
$ws.Range("AX7").Value = 0.96
$ws.Range("AY7").Value = 224223560
$ws.Range("AZ7").Value = 214466
$ws.Range("BA7").Value = 52
$ws.Range("BB7").Value = 218
$ws.Range("BC7").Value = 206
$ws.Range("BD7").Value = 270
$ws.Range("BE7").Value = 2108
$ws.Range("B8").Value = 14
$ws.Range("E8").Value = 14
$ws.Range("F8").Value = 13
$ws.Range("K8").Value = 1510
$ws.Range("P8").Value = 9.630000000000001
$ws.Range("AD8").Value = 9236
$ws.Range("AE8").Value = 0.89
$ws.Range("AG8").Value = 0.74
$ws.Range("AH8").Value = 5448
$ws.Range("AL8").Value = 99645
$ws.Range("AM8").Value = 0.27
$ws.Range("AO8").Value = 0.98
$ws.Range("AP8").Value = 7265
$ws.Range("AX8").Value = 1.26
$ws.Range("AY8").Value = 224223560
$ws.Range("AZ8").Value = 85918
$ws.Range("BA8").Value = 52
$ws.Range("BB8").Value = 437
$ws.Range("BC8").Value = 817
$ws.Range("BD8").Value = 1098
$ws.Range("BE8").Value = 3559
$ws.Range("K9").Value = 1510
$ws.Range("P9").Value = 9.630000000000001
$ws.Range("Q9").Value = 3658
$ws.Range("Z9").Value = "3658"
$ws.Range("AA9").Value = "254 x 146 x 43"
$ws.Range("AC9").Value = 3891
$ws.Range("AD9").Value = 4322
$ws.Range("AE9").Value = 0.95
$ws.Range("AG9").Value = 0.7
$ws.Range("AH9").Value = 2728
$ws.Range("AJ9").Value = 3891
$ws.Range("AK9").Value = 274000
$ws.Range("AL9").Value = 181090
$ws.Range("AP9").Value = 3936
$ws.Range("AR9").Value = 1945
$ws.Range("AS9").Value = 21924
$ws.Range("AT9").Value = 0.3
$ws.Range("AW9").Value = 1903
$ws.Range("AX9").Value = 0.88
$ws.Range("AY9").Value = 224223560
$ws.Range("AZ9").Value = 181090
$ws.Range("BA9").Value = 52
$ws.Range("BB9").Value = 218
$ws.Range("BC9").Value = 206
$ws.Range("BD9").Value = 263
$ws.Range("BE9").Value = 1855
$ws.Range("B10").Value = 14
$ws.Range("E10").Value = 14
$ws.Range("F10").Value = 13
$ws.Range("K10").Value = 1510
$ws.Range("P10").Value = 9.630000000000001
$ws.Range("AA10").Value = "457 x 152 x 52"
$ws.Range("AC10").Value = 4729
$ws.Range("AD10").Value = 3536
$ws.Range("AE10").Value = 1.16
$ws.Range("AG10").Value = 0.5600000000000001
$ws.Range("AH10").Value = 2639
$ws.Range("AJ10").Value = 4729
$ws.Range("AK10").Value = 333000
$ws.Range("AL10").Value = 63811
$ws.Range("AM10").Value = 0.27
$ws.Range("AO10").Value = 0.98
$ws.Range("AP10").Value = 4652
$ws.Range("AR10").Value = 2364
$ws.Range("AS10").Value = 20888
$ws.Range("AT10").Value = 0.34
$ws.Range("AW10").Value = 2291
$ws.Range("AX10").Value = 0.97
$ws.Range("AY10").Value = 224223560
$ws.Range("AZ10").Value = 55021
$ws.Range("BA10").Value = 52
$ws.Range("BB10").Value = 437
$ws.Range("BC10").Value = 817
$ws.Range("BD10").Value = 994
$ws.Range("BE10").Value = 2021
$ws.Range("K11").Value = 1510
$ws.Range("P11").Value = 9.630000000000001
$ws.Range("Q11").Value = 2038
$ws.Range("Z11").Value = "2038"
$ws.Range("AA11").Value = "254 x 102 x 25"
$ws.Range("AC11").Value = 2272
$ws.Range("AD11").Value = 2254
$ws.Range("AE11").Value = 1
$ws.Range("AG11").Value = 0.66
$ws.Range("AH11").Value = 1506
$ws.Range("AJ11").Value = 2272
$ws.Range("AK11").Value = 160000
$ws.Range("AL11").Value = 105746
$ws.Range("AP11").Value = 2298
$ws.Range("AR11").Value = 1136
$ws.Range("AS11").Value = 4825
$ws.Range("AT11").Value = 0.49
$ws.Range("AV11").Value = 0.93
$ws.Range("AW11").Value = 1055
$ws.Range("AX11").Value = 0.7
$ws.Range("AY11").Value = 224223560
$ws.Range("AZ11").Value = 105746
$ws.Range("BD11").Value = 235
$ws.Range("BE11").Value = 1042
